# Update the dSF (column F) values on Sheet1 to reflect the repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    5  = 0
    6  = -2
    10 = 2
    11 = -4
    13 = 1
    17 = -1
    18 = -3
    21 = 0
    22 = 4
    28 = -2
    29 = -4
    33 = -1
    37 = -2
    38 = -3
    41 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
